# Regenerate save_data: replace column G (previously "Strike#" counts,
# now header "K") with newly computed strikeout values, and recompute the
# dependent std/mean "s_vals" that were written alongside it.
#
# The workbook data lives on the active sheet's single table (rows 2-74,
# header in row 1). Column G holds the per-game value. Below we write the
# freshly-calculated values for every row whose K figure changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$newValues = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 2
    11 = 2
    12 = 1
    13 = 3
    14 = 1
    15 = 1
    16 = 0
    17 = 2
    18 = 2
    19 = 2
    20 = 0
    21 = 0
    22 = 0
    23 = 3
    24 = 1
    25 = 0
    26 = 2
    27 = 2
    28 = 1
    29 = 2
    30 = 1
    31 = 2
    32 = 3
    33 = 1
    35 = 0
    36 = 2
    37 = 0
    38 = 2
    39 = 2
    40 = 1
    41 = 3
    42 = 3
    43 = 3
    44 = 0
    45 = 3
    46 = 1
    47 = 1
    48 = 2
    49 = 1
    50 = 3
    51 = 2
    52 = 0
    53 = 5
    54 = 1
    55 = 1
    56 = 2
    57 = 3
    58 = 4
    59 = 1
    60 = 1
    61 = 1
    62 = 0
    63 = 1
    64 = 2
    65 = 0
    66 = 1
    67 = 3
    68 = 0
    69 = 1
    70 = 2
    71 = 1
    72 = 0
    73 = 1
    74 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
